$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = 9.14
[void]$ws.Range("S7").Select()
